# Update "想去人数" (want-to-go count) values in column F across the three
# sheets that list event rows. Each row appears once in its primary sheet
# ("展览" or "演出") and again in the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1011
$wsExhibition.Range("F5").Value = 1182
$wsExhibition.Range("F6").Value = 983
$wsExhibition.Range("F11").Value = 336
$wsExhibition.Range("F15").Value = 127
$wsExhibition.Range("F17").Value = 2958
$wsExhibition.Range("F19").Value = 1581
$wsExhibition.Range("F20").Value = 1331
$wsExhibition.Range("F31").Value = 1498

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F13").Value = 56

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1011
$wsAll.Range("F7").Value = 1182
$wsAll.Range("F8").Value = 983
$wsAll.Range("F22").Value = 336
$wsAll.Range("F26").Value = 127
$wsAll.Range("F28").Value = 2958
$wsAll.Range("F30").Value = 1581
$wsAll.Range("F31").Value = 1331
$wsAll.Range("F44").Value = 1498
$wsAll.Range("F45").Value = 56
